$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: A1 -> "Time", B1 -> "Input"
$ws.Range("A1").Value = "Time"
$ws.Range("B1").Value = "Input"

# Update column B data values (column A left unchanged)
$ws.Range("B2").Value = 6
$ws.Range("B3").Value = 165
$ws.Range("B4").Value = 118
$ws.Range("B5").Value = 37

# Update active selection to B7 (matches sheetView selection in target file)
$ws.Range("B7").Select()
